$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has daily data through row 67 (date serial 45623). Extend it with
# five more days (rows 68-72, serials 45624-45628), copying the rest of the
# row's values/format from the last existing row since they stay constant.
$srcRow = 67
$lastDate = $ws.Cells.Item($srcRow, 1).Value2

for ($i = 1; $i -le 5; $i++) {
    $newRow = $srcRow + $i

    $src = $ws.Range("A" + $srcRow + ":J" + $srcRow)
    $dst = $ws.Range("A" + $newRow + ":J" + $newRow)
    $src.Copy($dst)

    $ws.Cells.Item($newRow, 1).Value = $lastDate + $i
}
